# Form the consolidated report: fill in the "Absent" (column H) values
# based on the "Real" (column E) attendance values for each date row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H13").Value = 0
